$d = $word.ActiveDocument

# Paragraph 3 ("2022年6月2日星期四") is followed today by the paragraph that used
# to read "中雨，今天是农历五月初四，明天就是端午节了。" (old paragraph 4).
# The edit inserts two brand-new paragraphs between them - a repeat of that
# "中雨..." sentence, followed by a new "2022年6月3日星期五" date line - and
# then rewrites the (now pushed-down) original paragraph's text to describe
# the next day instead.

$p3 = $d.Paragraphs.Item(3)

# Create a single new empty paragraph right after paragraph 3; we will fill
# it (and grow it into two paragraphs) via a raw-XML insert so the run
# splits match exactly how Word split them (digits typed without the IME
# vs. the eastAsia-hinted Chinese runs around them).
$p3.Range.InsertParagraphAfter()
$newRange = $d.Paragraphs.Item(4).Range

$insertedXml = @'
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r>
              <w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr>
              <w:t>中雨，今天是农历五月初四，明天就是端午节了。</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr>
              <w:t>2</w:t>
            </w:r>
            <w:r>
              <w:t>022</w:t>
            </w:r>
            <w:r>
              <w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr>
              <w:t>年6月3日星期五</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$newRange.InsertXML($insertedXml)

# The original "中雨，今天是农历五月初四..." paragraph has now been pushed down
# to become the last paragraph; update only its text to the new sentence
# about the 5th day of the 5th lunar month / Dragon Boat Festival.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.Find.Execute("中雨，今天是农历五月初四，明天就是端午节了。", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "中雨，今天是农历五月初五，中国传统端午节。", 2)
